$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 is the theta_se row; row 6 is the lambda_se row.
# Update the standard-error values (text strings formatted like "(0.NN)")
# for the specifications that have new SE estimates after rerunning.

$ws.Range("C4").Value = "(0.52)"
$ws.Range("D4").Value = "(0.41)"
$ws.Range("E4").Value = "(0.37)"
$ws.Range("F4").Value = "(0.54)"
$ws.Range("G4").Value = "(0.56)"

$ws.Range("C6").Value = "(0.41)"
$ws.Range("D6").Value = "(0.33)"
$ws.Range("E6").Value = "(0.29)"
$ws.Range("F6").Value = "(0.45)"
$ws.Range("G6").Value = "(0.44)"
